$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on D and E columns so that numeric-looking
# strings (e.g. "214.29") are stored as text, matching the original
# inlineStr cell type, not auto-converted to a float by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.871.92"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "1.627.35"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "214.29"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").Value = "28.47"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "0.0607"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "0.0901"
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "1.860.96"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").Value = "1.626.72"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "29.872.24"
$ws.Range("E16").Value = "  +0.73%  "
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "63.99"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("D19").Value = "239.86"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E21").Value = "  +0.58%  "
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "9.76"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").Value = "2.16"
$ws.Range("E24").Value = "  +2.49%  "
$ws.Range("D25").Value = "157.54"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "1.422.45"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("E37").Value = "  -4.99%  "
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "0.555"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "74.41"
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").Value = "0.827"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "1.767.65"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D49").Value = "48.61"
$ws.Range("E49").Value = "  -9.94%  "
$ws.Range("D50").Value = "90.87"
$ws.Range("E51").Value = "  +9.44%  "

# Restore default (unstyled) cell style so no stray style index is
# left referencing the Text number format.
$ws.Range("D2:E51").Style = "Normal"
